# Update to match VSE Export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab ("workload" -> "Workload")
$ws.Name = "Workload"

# 2. Columns J/K: the old "Window" column is renamed/expanded into
#    "Backup Window (hours)" and moved to J, while "Years in Scope"
#    (formerly J) moves over to K. Swap both header text and the sample
#    row's values together so the underlying data stays paired correctly.
$ws.Range("J1").Value = "Backup Window (hours)"
$ws.Range("K1").Value = "Years in Scope"
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 3

# Give the swapped J/K header + data cells their own explicit formatting
# (bold + centered for the header, centered for the data row), matching
# the VSE export's look for this pair of columns.
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("K2").HorizontalAlignment = -4108

# 3. Update the rest of the sample data row to match the VSE export sample.
$ws.Range("A2").Value = "Site_A"
$ws.Range("B2").Value = "Site_B"
$ws.Range("C2").Value = "testWorkload"
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = 600
$ws.Range("G2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("M2").Value = "yes"
$ws.Range("N2").Value = "Yes"
$ws.Range("O2").Value = "No"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 30
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 6
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 1
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0

# 4. A couple of blank helper cells below the swapped columns, also seen in
#    the export (formatted the same way, just without a value).
$ws.Range("J3").HorizontalAlignment = -4108
$ws.Range("K3").HorizontalAlignment = -4108
$ws.Range("J4").HorizontalAlignment = -4108
$ws.Range("K4").HorizontalAlignment = -4108

# 5. View/window bookkeeping to mirror the VSE export snapshot.
$ws.Range("K5").Select()
